$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.48297233333333
$ws.Range("H2").Value = 43.448917
$ws.Range("I2").Value = 0.2019336017030403
$ws.Range("J2").Value = 0.2019336017030403
$ws.Range("M2").Value = 14.48297233333333
$ws.Range("N2").Value = 43.448917
$ws.Range("O2").Value = 0.2019336017030403
$ws.Range("P2").Value = 0.2019336017030403
$ws.Range("Q2").Value = 209.7564876080988
$ws.Range("R2").Value = 1887.808388472889
$ws.Range("S2").Value = 0.04077717949676212
$ws.Range("T2").Value = 0.04077717949676212
$ws.Range("G3").Value = 14.48297233333333
$ws.Range("H3").Value = 43.448917
$ws.Range("I3").Value = 0.2019336017030403
$ws.Range("J3").Value = 0.2019336017030403
$ws.Range("O3").Value = 0.007144147385663391
$ws.Range("P3").Value = 0.00714414738566339
$ws.Range("Q3").Value = 7.42091088324689
$ws.Range("R3").Value = 66.78819794922201
$ws.Range("S3").Value = 0.001442643412684368
$ws.Range("T3").Value = 0.001442643412684368
$ws.Range("G4").Value = 14.48297233333333
$ws.Range("H4").Value = 43.448917
$ws.Range("I4").Value = 0.2019336017030403
$ws.Range("J4").Value = 0.2019336017030403
$ws.Range("M4").Value = 56.726097
$ws.Range("N4").Value = 170.178291
$ws.Range("O4").Value = 0.7909222509112964
$ws.Range("P4").Value = 0.7909222509112963
$ws.Range("Q4").Value = 821.5624934289831
$ws.Range("R4").Value = 7394.062440860847
$ws.Range("S4").Value = 0.1597137787935938
$ws.Range("T4").Value = 0.1597137787935938
$ws.Range("I5").Value = 0.007144147385663391
$ws.Range("J5").Value = 0.00714414738566339
$ws.Range("M5").Value = 14.48297233333333
$ws.Range("N5").Value = 43.448917
$ws.Range("O5").Value = 0.2019336017030403
$ws.Range("P5").Value = 0.2019336017030403
$ws.Range("Q5").Value = 7.42091088324689
$ws.Range("R5").Value = 66.78819794922201
$ws.Range("S5").Value = 0.001442643412684368
$ws.Range("T5").Value = 0.001442643412684368
$ws.Range("I6").Value = 0.007144147385663391
$ws.Range("J6").Value = 0.00714414738566339
$ws.Range("O6").Value = 0.007144147385663391
$ws.Range("P6").Value = 0.00714414738566339
$ws.Range("S6").Value = 0.00005103884186808106
$ws.Range("T6").Value = 0.00005103884186808105
$ws.Range("I7").Value = 0.007144147385663391
$ws.Range("J7").Value = 0.00714414738566339
$ws.Range("M7").Value = 56.726097
$ws.Range("N7").Value = 170.178291
$ws.Range("O7").Value = 0.7909222509112964
$ws.Range("P7").Value = 0.7909222509112963
$ws.Range("Q7").Value = 29.065809207034
$ws.Range("R7").Value = 261.592282863306
$ws.Range("S7").Value = 0.005650465131110943
$ws.Range("T7").Value = 0.005650465131110941
$ws.Range("G8").Value = 56.726097
$ws.Range("H8").Value = 170.178291
$ws.Range("I8").Value = 0.7909222509112964
$ws.Range("J8").Value = 0.7909222509112963
$ws.Range("M8").Value = 14.48297233333333
$ws.Range("N8").Value = 43.448917
$ws.Range("O8").Value = 0.2019336017030403
$ws.Range("P8").Value = 0.2019336017030403
$ws.Range("Q8").Value = 821.5624934289831
$ws.Range("R8").Value = 7394.062440860847
$ws.Range("S8").Value = 0.1597137787935938
$ws.Range("T8").Value = 0.1597137787935938
$ws.Range("G9").Value = 56.726097
$ws.Range("H9").Value = 170.178291
$ws.Range("I9").Value = 0.7909222509112964
$ws.Range("J9").Value = 0.7909222509112963
$ws.Range("O9").Value = 0.007144147385663391
$ws.Range("P9").Value = 0.00714414738566339
$ws.Range("Q9").Value = 29.065809207034
$ws.Range("R9").Value = 261.592282863306
$ws.Range("S9").Value = 0.005650465131110943
$ws.Range("T9").Value = 0.005650465131110941
$ws.Range("G10").Value = 56.726097
$ws.Range("H10").Value = 170.178291
$ws.Range("I10").Value = 0.7909222509112964
$ws.Range("J10").Value = 0.7909222509112963
$ws.Range("M10").Value = 56.726097
$ws.Range("N10").Value = 170.178291
$ws.Range("O10").Value = 0.7909222509112964
$ws.Range("P10").Value = 0.7909222509112963
$ws.Range("Q10").Value = 3217.85008085341
$ws.Range("R10").Value = 28960.65072768068
$ws.Range("S10").Value = 0.6255580069865917
$ws.Range("T10").Value = 0.6255580069865916
